$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.102.59"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.892.89"
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.68"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5141"
$ws.Range("E7").Value = "  +1.91%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3754"
$ws.Range("E8").Value = "  +3.27%  "
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.23"
$ws.Range("E10").Value = "  +2.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9051"
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07637"
$ws.Range("E12").Value = "  +2.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.892.88"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "95.04"
$ws.Range("E14").Value = "  +2.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.264"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9990"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008490"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.44"
$ws.Range("E18").Value = "  +2.17%  "
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.126.37"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.067"
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.135.21"
$ws.Range("E22").Value = "  +1.49%  "
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.407"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.283"
$ws.Range("E25").Value = "  +10.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.59"
$ws.Range("E26").Value = "  -1.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.767"
$ws.Range("E27").Value = "  -1.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.04"
$ws.Range("E28").Value = "  +0.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.48"
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.948"
$ws.Range("E30").Value = "  +5.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.830"
$ws.Range("E31").Value = "  +3.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09181"
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05086"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.240"
$ws.Range("E34").Value = "  +7.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7800"
$ws.Range("E35").Value = "  +4.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.987"
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.291"
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("E38").Value = "  +4.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01998"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5600"
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.076"
$ws.Range("E41").Value = "  +0.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.102"
$ws.Range("E42").Value = "  +7.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.651"
$ws.Range("E43").Value = "  +2.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "117.71"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1509"
$ws.Range("E45").Value = "  +2.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4808"
$ws.Range("E46").Value = "  +2.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.19"
$ws.Range("E47").Value = "  +1.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9975"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.599"
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.51"
$ws.Range("E50").Value = "  +1.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.03"
$ws.Range("E51").Value = "  +1.38%  "
